# The workbook gains a new "New 200m transport /12, vel/depth in flux calc,
# dt = 1 d, j = 2" case, inserted as row 18 on Sheet1 (everything below it
# shifts down by one row).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a blank row at 18; Excel copies the formatting of row 17 (the row
# above) into it, which matches the styling the new case row needs.
$ws.Rows("18:18").Insert()

# Fill in the new case's label and data.
$ws.Range("B18").Value = "New 200m transport /12,  vel/depth in flux calc, dt = 1 d, j = 2"
$ws.Range("C18").Value = -0.3166
$ws.Range("D18").Value = 0.0319
$ws.Range("F18").Value = 0.0733
$ws.Range("J18").Value = -0.0159

# The other case rows that wrap their label text use a 32pt row height.
$ws.Rows("18:18").RowHeight = 32

# Leave the selection where the edit was made.
[void]$ws.Range("C19").Select()
